$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the weekly purchase analysis upgrade
$ws.Range("A2").Value = 43963
$ws.Range("B2").Value = 43971
$ws.Range("C2").Value = 1401065.028
$ws.Range("D2").Value = 119239.422
$ws.Range("E2").Value = 600189.9199999999
$ws.Range("F2").Value = 89747.91
$ws.Range("G2").Value = 830812.4570000001
$ws.Range("H2").Value = 244402.9349999999
$ws.Range("I2").Value = 0.4070136357725145
